{"js": "// The author's M2Doc field marker \"{m:userdoc 'zone1'}\" was stored as two\n// runs (\"{m\" and \":userdoc 'zone1'}\"). The new TokenIteratorFieldRewriterSplit\n// parser expects each token piece (\"{\", \"m\", \":userdoc 'zone1'\", \"}\") to live\n// in its own run, so split the paragraph's content into four runs with the\n// same text, keeping the trailing \"}\" run's whitespace-significant flag.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"{m:userdoc 'zone1'}\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find paragraph with text: \" + targetText);\n}\n\n// Rebuild the paragraph's content as four distinct runs via a raw OOXML\n// fragment (Office.js has no direct \"run\" object, so this is the reliable\n// way to force an exact run split without touching surrounding paragraphs).\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:r><w:t>{</w:t></w:r>' +\n              '<w:r><w:t>m</w:t></w:r>' +\n              \"<w:r><w:t>:userdoc 'zone1'</w:t></w:r>\" +\n              '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\ntarget.getRange(\"Content\").insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The author's M2Doc field marker \"{m:userdoc 'zone1'}\" was stored as two\n# runs (\"{m\" and \":userdoc 'zone1'}\"). The new TokenIteratorFieldRewriterSplit\n# parser expects each token piece (\"{\", \"m\", \":userdoc 'zone1'\", \"}\") to live\n# in its own run, so split the paragraph's content into four runs with the\n# same text, keeping the trailing \"}\" run's whitespace-significant flag.\n\n$d = $word.ActiveDocument\n\n$targetText = \"{m:userdoc 'zone1'}\"\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    if ($text -eq ($targetText + \"`r\")) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find paragraph with text: $targetText\"\n}\n\n# Range covering just the paragraph's content (excludes the trailing\n# paragraph-mark character) so the paragraph's own attributes are untouched.\n$start = $target.Range.Start\n$end = $target.Range.End - 1\n$contentRange = $d.Range($start, $end)\n\n# Rebuild the paragraph's content as four distinct runs via a raw OOXML\n# fragment (InsertXML replaces exactly the target range's contents) - this is\n# the reliable way to force an exact run split without merging back together.\n$xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:userdoc ''zone1''</w:t></w:r><w:r><w:t xml:space=\"preserve\">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$contentRange.InsertXML($xml)\n"}
